# Insert a new data row at row 38 (weekly update to the daily price log).
# Excel's Rows.Insert() shifts row 38 and everything below it down by one,
# bumping the sheet's used range from A1:T155 to A1:T156, and copies the
# formatting (incl. the date number-format on column D) of the row above
# into the freshly-inserted blank row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(38).Insert()

$ws.Range("A38").Value = 4
$ws.Range("B38").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C38").Value = 'Los Lagos'
$ws.Range("D38").Value = 45071
$ws.Range("E38").Value = 10
$ws.Range("F38").Value = 'Fruta'
$ws.Range("G38").Value = 100104
$ws.Range("H38").Value = 'Frutos de pepita'
$ws.Range("I38").Value = 100104003
$ws.Range("J38").Value = 'Membrillo'
$ws.Range("K38").Value = 'Champion'
$ws.Range("L38").Value = 'Primera'
$ws.Range("M38").Value = 200
$ws.Range("N38").Value = 13000
$ws.Range("O38").Value = 14000
$ws.Range("P38").Value = 13500
$ws.Range("Q38").Value = '$/caja 18 kilos empedrada'
$ws.Range("R38").Value = 'Región de O''Higgins'
$ws.Range("S38").Value = 750
$ws.Range("T38").Value = 18
